# Apply row 14-18 data refresh per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 111798760
$ws.Range("Q14").Value = 753108.8301749222
$ws.Range("R14").Value = 7091007.708399305
$ws.Range("S14").Value = 100
$ws.Range("AR14").Value = ""
$ws.Range("A15").Value = 111798795
$ws.Range("B15").Value = 81076
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5046
$ws.Range("F15").Value = "Grön jordtunga"
$ws.Range("G15").Value = "Microglossum viride"
$ws.Range("H15").Value = "(Pers.:Fr.) Gillet"
$ws.Range("AF15").Value = "'"
$ws.Range("A16").Value = 111798755
$ws.Range("B16").Value = 90709
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5448
$ws.Range("F16").Value = "Svartvit taggsvamp"
$ws.Range("G16").Value = "Phellodon connatus"
$ws.Range("H16").Value = "(Schultz) nom.prov"
$ws.Range("Q16").Value = 753030.7189070459
$ws.Range("R16").Value = 7090920.781295684
$ws.Range("S16").Value = 25
$ws.Range("AF16").Value = ""
$ws.Range("AI16").Value = ""
$ws.Range("A17").Value = 111961472
$ws.Range("B17").Value = 90709
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 5448
$ws.Range("F17").Value = "Svartvit taggsvamp"
$ws.Range("G17").Value = "Phellodon connatus"
$ws.Range("H17").Value = "(Schultz) nom.prov"
$ws.Range("I17").Value = "'"
$ws.Range("J17").Value = "'"
$ws.Range("AF17").Value = "'"
$ws.Range("AX17").Value = "Stefan Phalagorn Bergström, Annika  Carlberg , Andreas Estensen, Ola Elleström, Anne Järvinen, Emma Sewell, Thomas Strid"
$ws.Range("A18").Value = 111961716
$ws.Range("B18").Value = 81076
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 5046
$ws.Range("F18").Value = "Grön jordtunga"
$ws.Range("G18").Value = "Microglossum viride"
$ws.Range("H18").Value = "(Pers.:Fr.) Gillet"
$ws.Range("I18").Value = "'2"
$ws.Range("J18").Value = "mycel"
$ws.Range("AF18").Value = "mikroskoperad"
